$d = $word.ActiveDocument

# --- 1. Remove the bookmark that wrapped "UI_RESULTS" (signal-slot mapping sentence) ---
$bm = $d.Bookmarks.Item("__DdeLink__212_1482927881")
$bm.Delete()

# --- 2. Add a new "See also" heading + the PyQt5 signals/slots reference link ---
# Locate the (empty) TextBody paragraph that follows the "...configurational part
# of the interface between UI elements and signals." paragraph.
$idx = 0
$targetIdx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -match "such that the enum is the configurational part of the interface between UI elements and signals\.") {
        $targetIdx = $idx + 1
    }
}

$target = $d.Paragraphs.Item($targetIdx)

# Insert a new paragraph before it and make it the "See also" heading.
$target.Range.InsertParagraphBefore()
$heading = $d.Paragraphs.Item($targetIdx)
$heading.Range.Text = "See also"
$heading.Style = "Heading"

# Fill the (still empty) TextBody paragraph with the reference URL.
$urlPara = $d.Paragraphs.Item($targetIdx + 1)
$urlPara.Range.InsertAfter("https://www.riverbankcomputing.com/static/Docs/PyQt5/signals_slots.html")
